$wb = $excel.ActiveWorkbook

# --- 1. Add a row to "Test Cases" describing the new test sheet ---
$tc = $wb.Worksheets.Item("Test Cases")
$tc.Range("A3").Value = "GoogleSearch1Test"
$tc.Range("B3").Value = "xxxxxxxxx"
$tc.Range("C3").Value = "Y"
$tc.Range("D3").Value = "PASS"
$tc.Range("D3").Style = $tc.Range("D2").Style

# --- 2. Create the new worksheet "GoogleSearch1Test" after "GoogleSearchTest" ---
$src = $wb.Worksheets.Item("GoogleSearchTest")
$newSheet = $wb.Worksheets.Add($null, $src)
$newSheet.Name = "GoogleSearch1Test"

$newSheet.Range("A1").Value = "Data1"
$newSheet.Range("B1").Value = "Data2"
$newSheet.Range("C1").Value = "Data3"
$newSheet.Range("D1").Value = "Data4"
$newSheet.Range("E1").Value = "Runmode"
$newSheet.Range("F1").Value = "Results"
$newSheet.Range("A1:F1").Style = $src.Range("A1:F1").Style

$newSheet.Range("A2").Value = "Chrome"
$newSheet.Range("B2").Value = "http://google.com"
$newSheet.Range("C2").Value = "hello world"
$newSheet.Range("D2").Value = "Google"
$newSheet.Range("E2").Value = "Y"
$newSheet.Range("F2").Value = "PASS"
$newSheet.Range("F2").Style = $src.Range("F2").Style

$newSheet.Range("A3").Value = "Chrome"
$newSheet.Range("B3").Value = "http://google.com"
$newSheet.Range("C3").Value = "hello world"
$newSheet.Range("D3").Value = "not Google"
$newSheet.Range("E3").Value = "Y"
$newSheet.Range("F3").Value = "PASS"

$newSheet.Hyperlinks.Add($newSheet.Range("B2"), "http://google.com")
$newSheet.Hyperlinks.Add($newSheet.Range("B3"), "http://google.com")
$newSheet.Range("B2").Style = $src.Range("B2").Style
$newSheet.Range("B3").Style = $src.Range("B2").Style

$newSheet.Columns.Item(6).ColumnWidth = 6.76953125

# --- 3. Selections / active tab ---
$tc.Range("A3").Select()
$src.Range("A1:F3").Select()
$newSheet.Range("A1:F3").Select()
$newSheet.Activate()

$wb.Worksheets.Item(1).Select()
